$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the crypto table (rows 2-51) per the latest scrape.
# Numeric-looking price strings are prefixed with a literal leading
# apostrophe so Excel stores them as text (preserving formatting like
# trailing zeros / decimal grouping) instead of auto-converting them
# to numbers.

$ws.Range("D2").Value = "34.626.20"
$ws.Range("E2").Value = "  +1.78%  "

$ws.Range("D3").Value = "1.789.14"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").Value = "'224.93"
$ws.Range("E5").Value = "  -0.17%  "

$ws.Range("D6").Value = "'0.556"
$ws.Range("E6").Value = "  -0.33%  "

$ws.Range("E7").Value = "  +0.25%  "

$ws.Range("D8").Value = "'32.65"

$ws.Range("D9").Value = "'0.284"
$ws.Range("E9").Value = "  +1.51%  "

$ws.Range("E10").Value = "  +0.72%  "

$ws.Range("E11").Value = "  +1.36%  "

$ws.Range("D12").Value = "2.047.42"
$ws.Range("E12").Value = "  +0.34%  "

$ws.Range("D13").Value = "'11.08"
$ws.Range("E13").Value = "  +10.43%  "

$ws.Range("D14").Value = "1.794.52"
$ws.Range("E14").Value = "  +0.28%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.633"
$ws.Range("E15").Value = "  +0.54%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "34.590.87"
$ws.Range("E16").Value = "  +1.87%  "

$ws.Range("D17").Value = "'4.29"
$ws.Range("E17").Value = "  +1.88%  "

$ws.Range("D18").Value = "'68.95"
$ws.Range("E18").Value = "  +0.30%  "

$ws.Range("D19").Value = "'254.07"
$ws.Range("E19").Value = "  +0.62%  "

$ws.Range("D20").Value = "0.0₃0766"
$ws.Range("E20").Value = "  +3.10%  "

$ws.Range("E21").Value = "  +0.15%  "

$ws.Range("D22").Value = "'10.40"
$ws.Range("E22").Value = "  +0.71%  "

$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("E24").Value = "  -1.37%  "

$ws.Range("D25").Value = "'159.19"
$ws.Range("E25").Value = "  +0.45%  "

$ws.Range("D26").Value = "'16.40"
$ws.Range("E26").Value = "  -0.74%  "

$ws.Range("D27").Value = "'7.08"
$ws.Range("E27").Value = "  +1.57%  "

$ws.Range("E28").Value = "  -0.72%  "

$ws.Range("E29").Value = "  +0.21%  "

$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "'3.77"
$ws.Range("E30").Value = "  -0.89%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.0516"
$ws.Range("E31").Value = "  +0.31%  "

$ws.Range("E32").Value = "  -0.33%  "

$ws.Range("E33").Value = "  +0.67%  "

$ws.Range("D34").Value = "'1.84"
$ws.Range("E34").Value = "  +3.08%  "

$ws.Range("D35").Value = "1.441.35"
$ws.Range("E35").Value = "  -3.77%  "

$ws.Range("D36").Value = "'1.06"
$ws.Range("E36").Value = "  -0.17%  "

$ws.Range("D37").Value = "'0.0189"
$ws.Range("E37").Value = "  +2.00%  "

$ws.Range("D38").Value = "'0.627"
$ws.Range("E38").Value = "  -1.02%  "

$ws.Range("D39").Value = "'82.86"
$ws.Range("E39").Value = "  -0.59%  "

$ws.Range("E40").Value = "  +4.43%  "

$ws.Range("D41").Value = "'2.35"
$ws.Range("E41").Value = "  -0.11%  "

$ws.Range("E42").Value = "  +0.84%  "

$ws.Range("D43").Value = "'2.12"
$ws.Range("E43").Value = "  +1.10%  "

$ws.Range("D44").Value = "'0.0504"
$ws.Range("E44").Value = "  -1.22%  "

$ws.Range("D45").Value = "'5.91"
$ws.Range("E45").Value = "  +2.04%  "

$ws.Range("E46").Value = "  -1.46%  "

$ws.Range("D47").Value = "1.947.53"
$ws.Range("E47").Value = "  +0.57%  "

$ws.Range("E48").Value = "  +0.28%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'11.91"
$ws.Range("E49").Value = "  -0.51%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'102.78"
$ws.Range("E50").Value = "  +4.70%  "

$ws.Range("D51").Value = "0.0₆0122"
$ws.Range("E51").Value = "  +4.10%  "
